# ---------------------------------------------------------------------------
# feat: add 2022-Q3 data
#
# 1. Insert a new "2022-Q3" worksheet (with its fund-holdings detail table)
#    right after the "总计" summary sheet and before "2021-Q4".
# 2. Prepend a new summary row for 2022-Q3 on the "总计" sheet and renumber
#    the existing index column.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Part 1: create the "2022-Q3" worksheet with its fund table
# ---------------------------------------------------------------------------

$beforeSheet = $wb.Worksheets.Item(2)
$newSheet = $wb.Worksheets.Add($beforeSheet)
$newSheet.Name = "2022-Q3"

# Re-fetch a live handle on the (still-present) "2021-Q4" sheet -- the
# object returned by Add() can alias the variable used for "Before", so
# look sheets back up by name/index instead of reusing pre-Add references.
$origSheet = $wb.Worksheets.Item("2021-Q4")

# Clone the header-row and index-column formatting (bold font + border)
# used by every quarterly detail sheet.
$origSheet.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)
$origSheet.Range("A2").Copy()
$newSheet.Range("A2:A20").PasteSpecial(-4122)

$ws = $newSheet
$ws.Range("B1").Value = "基金代码"
$ws.Range("C1").Value = "基金名称"
$ws.Range("D1").Value = "基金规模"
$ws.Range("E1").Value = "股票总仓位"
$ws.Range("F1").Value = "仓位占比"
$ws.Range("G1").Value = "持有市值(亿元)"
$ws.Range("H1").Value = "仓位排名"

# Columns B, D, E, F, G hold text that looks numeric (fund codes / percents
# formatted as strings) -- force Text format so COM keeps them as strings
# instead of silently coercing to numbers.
$ws.Range("B2:B20").NumberFormat = "@"
$ws.Range("D2:G20").NumberFormat = "@"

$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "720001"
$ws.Range("C2").Value = "财通价值动量混合"
$ws.Range("D2").Value = "38.35"
$ws.Range("E2").Value = "79.13"
$ws.Range("F2").Value = "5.96"
$ws.Range("G2").Value = "2.2857"
$ws.Range("H2").Value = 6
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "001480"
$ws.Range("C3").Value = "财通成长优选混合"
$ws.Range("D3").Value = "20.31"
$ws.Range("E3").Value = "91.20"
$ws.Range("F3").Value = "5.77"
$ws.Range("G3").Value = "1.1719"
$ws.Range("H3").Value = 7
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "014915"
$ws.Range("C4").Value = "财通匠心优选一年持有期混合A"
$ws.Range("D4").Value = "5.65"
$ws.Range("E4").Value = "81.89"
$ws.Range("F4").Value = "6.33"
$ws.Range("G4").Value = "0.3576"
$ws.Range("H4").Value = 6
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "009970"
$ws.Range("C5").Value = "财通内需增长12个月定期开放混合"
$ws.Range("D5").Value = "9.38"
$ws.Range("E5").Value = "56.38"
$ws.Range("F5").Value = "3.42"
$ws.Range("G5").Value = "0.3208"
$ws.Range("H5").Value = 2
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = "008983"
$ws.Range("C6").Value = "财通科技创新混合A"
$ws.Range("D6").Value = "2.95"
$ws.Range("E6").Value = "94.28"
$ws.Range("F6").Value = "6.51"
$ws.Range("G6").Value = "0.1920"
$ws.Range("H6").Value = 3
$ws.Range("A7").Value = 5
$ws.Range("B7").Value = "501046"
$ws.Range("C7").Value = "财通多策略福鑫定期开放灵活配置混合"
$ws.Range("D7").Value = "2.82"
$ws.Range("E7").Value = "85.55"
$ws.Range("F7").Value = "6.81"
$ws.Range("G7").Value = "0.1920"
$ws.Range("H7").Value = 5
$ws.Range("A8").Value = 6
$ws.Range("B8").Value = "009062"
$ws.Range("C8").Value = "财通智慧成长混合A"
$ws.Range("D8").Value = "2.17"
$ws.Range("E8").Value = "84.78"
$ws.Range("F8").Value = "6.63"
$ws.Range("G8").Value = "0.1439"
$ws.Range("H8").Value = 5
$ws.Range("A9").Value = 7
$ws.Range("B9").Value = "009063"
$ws.Range("C9").Value = "财通智慧成长混合C"
$ws.Range("D9").Value = "1.50"
$ws.Range("E9").Value = "84.78"
$ws.Range("F9").Value = "6.63"
$ws.Range("G9").Value = "0.0994"
$ws.Range("H9").Value = 5
$ws.Range("A10").Value = 8
$ws.Range("B10").Value = "000017"
$ws.Range("C10").Value = "财通可持续发展主题混合"
$ws.Range("D10").Value = "1.95"
$ws.Range("E10").Value = "90.33"
$ws.Range("F10").Value = "4.03"
$ws.Range("G10").Value = "0.0786"
$ws.Range("H10").Value = 5
$ws.Range("A11").Value = 9
$ws.Range("B11").Value = "008984"
$ws.Range("C11").Value = "财通科技创新混合C"
$ws.Range("D11").Value = "1.16"
$ws.Range("E11").Value = "94.28"
$ws.Range("F11").Value = "6.51"
$ws.Range("G11").Value = "0.0755"
$ws.Range("H11").Value = 3
$ws.Range("A12").Value = 10
$ws.Range("B12").Value = "013238"
$ws.Range("C12").Value = "财通均衡一年持有期混合A"
$ws.Range("D12").Value = "1.81"
$ws.Range("E12").Value = "86.47"
$ws.Range("F12").Value = "4.01"
$ws.Range("G12").Value = "0.0726"
$ws.Range("H12").Value = 5
$ws.Range("A13").Value = 11
$ws.Range("B13").Value = "501026"
$ws.Range("C13").Value = "财通多策略福享混合（LOF）"
$ws.Range("D13").Value = "1.43"
$ws.Range("E13").Value = "88.50"
$ws.Range("F13").Value = "4.02"
$ws.Range("G13").Value = "0.0575"
$ws.Range("H13").Value = 5
$ws.Range("A14").Value = 12
$ws.Range("B14").Value = "014916"
$ws.Range("C14").Value = "财通匠心优选一年持有期混合C"
$ws.Range("D14").Value = "0.61"
$ws.Range("E14").Value = "81.89"
$ws.Range("F14").Value = "6.33"
$ws.Range("G14").Value = "0.0386"
$ws.Range("H14").Value = 6
$ws.Range("A15").Value = 13
$ws.Range("B15").Value = "000916"
$ws.Range("C15").Value = "前海开源股息率100强等权重股票"
$ws.Range("D15").Value = "2.44"
$ws.Range("E15").Value = "94.04"
$ws.Range("F15").Value = "1.33"
$ws.Range("G15").Value = "0.0325"
$ws.Range("H15").Value = 8
$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "006968"
$ws.Range("C16").Value = "财通行业龙头精选混合C"
$ws.Range("D16").Value = "0.07"
$ws.Range("E16").Value = "90.94"
$ws.Range("F16").Value = "4.06"
$ws.Range("G16").Value = "0.0028"
$ws.Range("H16").Value = 6
$ws.Range("A17").Value = 15
$ws.Range("B17").Value = "006967"
$ws.Range("C17").Value = "财通行业龙头精选混合A"
$ws.Range("D17").Value = "0.05"
$ws.Range("E17").Value = "90.94"
$ws.Range("F17").Value = "4.06"
$ws.Range("G17").Value = "0.0020"
$ws.Range("H17").Value = 6
$ws.Range("A18").Value = 16
$ws.Range("B18").Value = "013239"
$ws.Range("C18").Value = "财通均衡一年持有期混合C"
$ws.Range("D18").Value = "0.04"
$ws.Range("E18").Value = "86.47"
$ws.Range("F18").Value = "4.01"
$ws.Range("G18").Value = "0.0016"
$ws.Range("H18").Value = 5
$ws.Range("A19").Value = 17
$ws.Range("B19").Value = "005126"
$ws.Range("C19").Value = "银河量化稳进混合"
$ws.Range("D19").Value = "0.13"
$ws.Range("E19").Value = "55.69"
$ws.Range("F19").Value = "1.14"
$ws.Range("G19").Value = "0.0015"
$ws.Range("H19").Value = 3
$ws.Range("A20").Value = 18
$ws.Range("B20").Value = "001849"
$ws.Range("C20").Value = "前海开源强势共识100强等权重股票"
$ws.Range("D20").Value = "0.11"
$ws.Range("E20").Value = "91.68"
$ws.Range("F20").Value = "1.04"
$ws.Range("G20").Value = "0.0011"
$ws.Range("H20").Value = 4

# ---------------------------------------------------------------------------
# Part 2: update the "总计" summary sheet
# ---------------------------------------------------------------------------

$total = $wb.Worksheets.Item("总计")

# Insert a new row 2 (pushing the three existing quarters down) and copy
# the row-3 formatting onto it so the index column keeps its style.
$total.Rows.Item(2).Insert()
$total.Range("A3:D3").Copy()
$total.Range("A2:D2").PasteSpecial(-4122)

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q3"
$total.Range("C2").Value = 19
$total.Range("D2").Value = 5.13

$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
$total.Range("A5").Value = 3

# ---------------------------------------------------------------------------
# Restore the originally-active sheet (adding a worksheet activates it).
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("2020-Q4").Select()

Write-Output "2022-Q3 sheet + summary row added"
